# Update the course-number prefix "107" -> "108" in every slide title
# (e.g. "107 Final Project Review" -> "108 Final Project Review",
#  "107.1 Data Preparation" -> "108.1 Data Preparation", etc.)
#
# The text lives in the title placeholder (always the first shape on the
# slides that carry it), as a single text run. We walk every slide/shape,
# and for any text frame whose text begins with "107" we rewrite the
# leading "107" to "108", leaving the remainder (and all run formatting,
# since we only touch the text) untouched.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shape = $s.Shapes.Item($j)

        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $text = $tr.Text

        if ($text -match '^107(\.\d+)?( |$)') {
            $tr.Text = $text -replace '^107', '108'
        }
    }
}
